# edit.ps1 - applies the WIP docx-generator template changes:
#  1. Table cell paragraph "ที่ สพจ..." -> split "สพจ" into its own
#     spell-checked run (proofErr spellStart/spellEnd) and move the
#     "_GoBack" bookmark here (from its old location further down).
#  2. Paragraph holding "{#attachment}" loses the "_GoBack" bookmark
#     (it now lives on paragraph 1, see above).
#  3. Footer paragraph: "{contact_name}" and "{contact_phone}" merge
#     fields get split into "{" + spell-checked placeholder + "}" runs.

$d = $word.ActiveDocument

function Replace-ParagraphXml($range, [string]$fragment) {
    $pkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
           "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
           "<pkg:xmlData>" +
           "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'>" +
           "<w:body>" + $fragment + "</w:body></w:document>" +
           "</pkg:xmlData></pkg:part></pkg:package>"
    $range.InsertXML($pkg)
}

# --- 1 & 2: main document body -------------------------------------------------

$para1Xml = '<w:p w14:paraId="75837872" w14:textId="3BA5AC5F" w:rsidR="00AE1FFA" w:rsidRPr="00D5432C" w:rsidRDefault="00AE1FFA" w:rsidP="002077A5"><w:pPr><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="00D5432C"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/><w:cs/></w:rPr><w:t xml:space="preserve">ที่ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00D5432C"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/><w:cs/></w:rPr><w:t>สพจ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00D5432C"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="00416E23"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/></w:rPr><w:t>{number}</w:t></w:r></w:p>'
$paraAttachXml = '<w:p w14:paraId="27E717BB" w14:textId="6D71C32E" w:rsidR="00AE1FFA" w:rsidRDefault="00AE1FFA" w:rsidP="00AE1FFA"><w:pPr><w:pStyle w:val="textbox"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00D5432C"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:cs/></w:rPr><w:t>เรียน</w:t></w:r><w:r w:rsidRPr="00D5432C"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00D5432C"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:cs/></w:rPr><w:tab/></w:r><w:r w:rsidR="005E5827"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>{to}</w:t></w:r><w:r w:rsidR="00092E5A"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>{#attachment}</w:t></w:r></w:p>'

# Paragraph 1: "ที่ สพจ. {number}" table cell
$target1 = $d.Content
$target1.Find.Execute("{number}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para1 = $target1.Paragraphs(1).Range
Replace-ParagraphXml $para1 $para1Xml

# Paragraph containing "{#attachment}"
$targetA = $d.Content
$targetA.Find.Execute("{#attachment}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraA = $targetA.Paragraphs(1).Range
Replace-ParagraphXml $paraA $paraAttachXml

# --- 3: footer -------------------------------------------------------------

$footerXml = '<w:p w14:paraId="494AEA31" w14:textId="16791B23" w:rsidR="00E121EA" w:rsidRPr="003B3D66" w:rsidRDefault="00E121EA" w:rsidP="00244F75"><w:pPr><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:spacing w:val="5"/><w:szCs w:val="24"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="003A66A6"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New" w:hint="cs"/><w:szCs w:val="24"/><w:cs/></w:rPr><w:t>นิ</w:t></w:r><w:r w:rsidRPr="003A66A6"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/><w:cs/></w:rPr><w:t>สิตผู้ประสานงาน</w:t></w:r><w:r w:rsidRPr="0045270E"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="0045270E"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/><w:cs/></w:rPr><w:tab/></w:r><w:r w:rsidR="005E5827"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="005E5827"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/></w:rPr><w:t>contact_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="005E5827"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/></w:rPr><w:t>}</w:t></w:r><w:r w:rsidR="00244F75" w:rsidRPr="00244F75"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/><w:cs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00244F75"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New" w:hint="cs"/><w:szCs w:val="24"/><w:cs/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidR="00244F75" w:rsidRPr="00244F75"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/><w:cs/></w:rPr><w:t xml:space="preserve">หมายเลขโทรศัพท์ </w:t></w:r><w:r w:rsidR="005E5827"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="005E5827"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/></w:rPr><w:t>contact_phone</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="005E5827"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New"/><w:szCs w:val="24"/></w:rPr><w:t>}</w:t></w:r><w:r w:rsidR="00244F75"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:eastAsia="TH SarabunPSK" w:hAnsi="TH Sarabun New" w:cs="TH Sarabun New" w:hint="cs"/><w:szCs w:val="24"/><w:cs/></w:rPr><w:t>)</w:t></w:r></w:p>'
$ftr = $d.Sections(1).Footers(1)
$footerPara = $ftr.Range.Paragraphs(1).Range
Replace-ParagraphXml $footerPara $footerXml

Write-Host "done"
